$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 5
    3  = 2
    4  = 6
    5  = 4
    6  = 1
    8  = 6
    9  = 4
    10 = 3
    11 = 6
    12 = 5
    13 = 3
    14 = 5
    16 = 1
    17 = 3
    18 = 3
    19 = 8
    20 = 7
    21 = 2
    22 = 1
    23 = 6
    24 = 5
    25 = 5
    26 = 6
    28 = 3
    29 = 5
    30 = 4
    32 = 4
    33 = 6
    34 = 4
    35 = 1
    36 = 7
    37 = 3
    38 = 8
    39 = 5
    40 = 7
    41 = 3
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 10).Value = $values[$row]
}
